$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 38000
$ws.Range("J75").Value = 38000
$ws.Range("L75").Value = 38000
$ws.Range("N75").Value = -39872
$ws.Range("H78").Value = 38000
$ws.Range("J78").Value = 38000
$ws.Range("L78").Value = 114000
$ws.Range("N78").Value = -123360
$ws.Range("H114").Value = 36666.668
$ws.Range("J114").Value = 36666.668
$ws.Range("L114").Value = 36666.668
$ws.Range("N114").Value = -45344.668
$ws.Range("H126").Value = 54245
$ws.Range("J126").Value = 54245
$ws.Range("L126").Value = 54245
$ws.Range("N126").Value = -64125
$ws.Range("H134").Value = 39136.332
$ws.Range("J134").Value = 39136.332
$ws.Range("L134").Value = 39136.332
$ws.Range("N134").Value = -49276.332
$ws.Range("H141").Value = 2690.3572
$ws.Range("I141").Value = 2690.3572
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8071.071599999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2891.071599999999
$ws.Range("N141").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 39478
$ws.Range("J103").Value = 39478
$ws.Range("L103").Value = 39478
$ws.Range("N103").Value = -41822
$ws.Range("H122").Value = 2329.7646
$ws.Range("I122").Value = 2088.25
$ws.Range("K122").Value = 6264.75
$ws.Range("M122").Value = -3814.75
$ws.Range("H129").Value = 43996.855
$ws.Range("I129").Value = 40000
$ws.Range("J129").Value = 44663
$ws.Range("K129").Value = 40000
$ws.Range("L129").Value = 44663
$ws.Range("M129").Value = -35000
$ws.Range("N129").Value = -54663
$ws.Range("H133").Value = 55800
$ws.Range("J133").Value = 55800
$ws.Range("L133").Value = 55800
$ws.Range("N133").Value = -60860

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 983.1667
$ws.Range("I22").Value = 979.8
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 979.8
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -806.8
$ws.Range("N22").Value = -1346
$ws.Range("H122").Value = 42024.547
$ws.Range("J122").Value = 42024.547
$ws.Range("L122").Value = 42024.547
$ws.Range("N122").Value = -51824.547
$ws.Range("H125").Value = 53225
$ws.Range("J125").Value = 53225
$ws.Range("L125").Value = 53225
$ws.Range("N125").Value = -63065
$ws.Range("H132").Value = 55493.332
$ws.Range("J132").Value = 55493.332
$ws.Range("L132").Value = 55493.332
$ws.Range("N132").Value = -65613.33199999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 16582.5
$ws.Range("J41").Value = 16582.5
$ws.Range("L41").Value = 16582.5
$ws.Range("N41").Value = -17438.5
$ws.Range("H68").Value = 24874.666
$ws.Range("J68").Value = 24874.666
$ws.Range("L68").Value = 24874.666
$ws.Range("N68").Value = -26372.666
$ws.Range("H71").Value = 24874.666
$ws.Range("J71").Value = 24874.666
$ws.Range("L71").Value = 74623.99800000001
$ws.Range("N71").Value = -82111.99800000001
$ws.Range("H123").Value = 43750
$ws.Range("J123").Value = 43750
$ws.Range("L123").Value = 43750
$ws.Range("N123").Value = -53550

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 816.6667
$ws.Range("I5").Value = 633.3333
$ws.Range("K5").Value = 1899.9999
$ws.Range("M5").Value = -1787.9999
$ws.Range("H135").Value = 816.6667
$ws.Range("I135").Value = 633.3333
$ws.Range("K135").Value = 5699.9997
$ws.Range("M135").Value = -3164.9997
$ws.Range("H137").Value = 1426.625
$ws.Range("I137").Value = 1278.1818
$ws.Range("J137").Value = 1753.2
$ws.Range("K137").Value = 3834.5454
$ws.Range("L137").Value = 5259.6
$ws.Range("M137").Value = 1265.4546
$ws.Range("N137").Value = -15459.6

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 19071.666
$ws.Range("J62").Value = 19071.666
$ws.Range("L62").Value = 19071.666
$ws.Range("N62").Value = -20443.666
$ws.Range("H65").Value = 19071.666
$ws.Range("J65").Value = 19071.666
$ws.Range("L65").Value = 57214.99800000001
$ws.Range("N65").Value = -64078.99800000001
$ws.Range("H113").Value = 1705.7894
$ws.Range("I113").Value = 1682.9412
$ws.Range("K113").Value = 1682.9412
$ws.Range("M113").Value = 487.0588
$ws.Range("H122").Value = 3099.5557
$ws.Range("I122").Value = 3399.3333
$ws.Range("K122").Value = 10197.9999
$ws.Range("M122").Value = -7747.999899999999
$ws.Range("H124").Value = 54880
$ws.Range("J124").Value = 54880
$ws.Range("L124").Value = 54880
$ws.Range("N124").Value = -64700
$ws.Range("H127").Value = 48740
$ws.Range("J127").Value = 48740
$ws.Range("L127").Value = 48740
$ws.Range("N127").Value = -58660
$ws.Range("H130").Value = 55521.25
$ws.Range("J130").Value = 55521.25
$ws.Range("L130").Value = 55521.25
$ws.Range("N130").Value = -65561.25
$ws.Range("H133").Value = 39000
$ws.Range("J133").Value = 39000
$ws.Range("L133").Value = 39000
$ws.Range("N133").Value = -49120

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 32450
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 32450
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 32450
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -33948
$ws.Range("H66").Value = 32450
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 32450
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 97350
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -104838
$ws.Range("H82").Value = 1550.579
$ws.Range("I82").Value = 1119.9
$ws.Range("J82").Value = 2029.1111
$ws.Range("K82").Value = 1119.9
$ws.Range("L82").Value = 2029.1111
$ws.Range("M82").Value = -758.9000000000001
$ws.Range("N82").Value = -2751.1111
$ws.Range("H85").Value = 1550.579
$ws.Range("I85").Value = 1119.9
$ws.Range("J85").Value = 2029.1111
$ws.Range("K85").Value = 1119.9
$ws.Range("L85").Value = 2029.1111
$ws.Range("M85").Value = 128.0999999999999
$ws.Range("N85").Value = -4525.1111
$ws.Range("H93").Value = 2463.125
$ws.Range("I93").Value = 2264.182
$ws.Range("J93").Value = 2900.8
$ws.Range("K93").Value = 2264.182
$ws.Range("L93").Value = 2900.8
$ws.Range("M93").Value = -1016.182
$ws.Range("N93").Value = -5396.8
$ws.Range("H96").Value = 11566.167
$ws.Range("J96").Value = 11566.167
$ws.Range("L96").Value = 11566.167
$ws.Range("N96").Value = -17058.167
$ws.Range("H108").Value = 24326.666
$ws.Range("J108").Value = 24326.666
$ws.Range("L108").Value = 24326.666
$ws.Range("N108").Value = -32006.666
$ws.Range("H129").Value = 33247.6
$ws.Range("J129").Value = 33247.6
$ws.Range("L129").Value = 33247.6
$ws.Range("N129").Value = -43247.6

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H62").Value = 3588.889
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3588.889
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3588.889
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4836.889
$ws.Range("H65").Value = 3588.889
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3588.889
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 17944.445
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -24184.445
$ws.Range("H122").Value = 3533.625
$ws.Range("I122").Value = 2603.8
$ws.Range("J122").Value = 5083.3335
$ws.Range("K122").Value = 7811.400000000001
$ws.Range("L122").Value = 15250.0005
$ws.Range("M122").Value = -5361.400000000001
$ws.Range("N122").Value = -20150.0005
$ws.Range("H125").Value = 49418.75
$ws.Range("J125").Value = 49418.75
$ws.Range("L125").Value = 49418.75
$ws.Range("N125").Value = -59258.75
$ws.Range("H129").Value = 25845.666
$ws.Range("J129").Value = 25845.666
$ws.Range("L129").Value = 25845.666
$ws.Range("N129").Value = -35845.666
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
